$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.073.25'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -5.30%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.826.63'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.99%  '

$ws.Range('E4').Value = '  -0.90%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '327.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.28%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.88%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4615'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3853'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.95%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.89%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07829'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.28%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9588'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.41%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.86'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.53%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.675'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.38%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.871'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.54%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.744.05'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -11.32%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06843'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.08%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.003'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.98%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '86.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.18%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000009917'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.53%  '

$ws.Range('E20').Value = '  -4.22%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.00%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.095.97'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.30%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.316'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.61%  '

$ws.Range('E24').Value = '  -7.05%  '

$ws.Range('E25').Value = '  -1.60%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.005.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -9.16%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '152.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.18%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.60%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.712'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -13.27%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.971'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.53%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '116.57'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.13%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09272'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.38%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9360'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.53%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.250'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.88%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.427'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.20%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.309'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.29%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05967'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.82%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02142'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.38%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.144'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.71%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('D40').Style = 'Normal'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.557'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.15%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5592'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.38%  '

$ws.Range('E43').Value = '  -6.80%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1765'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.70%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.221'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.12%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.244'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -9.66%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.63'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.33%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5248'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.59%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06998'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.28%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.827'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '112.21'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.60%  '
